# edit.ps1
# Applies the "ops per 1 ms, not per 100 / in superscript" commit:
#   1. Remove the paragraph "Also, for my script I standardized it to 100
#      milliseconds instead of 1 millisecond." along with the blank paragraph
#      that immediately follows it.
#   2. Reword "My throughput improved from 73 operations in 100 milliseconds
#      to 88 operations." -> "My throughput improved from .73 operations per
#      millisecond to .88 operations."
#   3. Reword "operations per 100 millisecs: 73" -> "avg operations per
#      millisec: .73"
#   4. Reword " operations per 100 millisecs: 88" -> " avg operations per
#      millisec: .88" and, on that same (last "Running superscript")
#      paragraph, drop the stray <w:color w:val="000000"/> from the
#      paragraph mark run-properties while adding the widowControl /
#      autoSpaceDE / autoSpaceDN / adjustRightInd paragraph properties that
#      match the rest of the test-output paragraphs in the doc.
#   5. Replace the final paragraph (a lone " " run in Monaco/color 000000)
#      with a truly empty paragraph that carries the same widowControl /
#      autoSpaceDE / autoSpaceDN / adjustRightInd paragraph properties.

$d = $word.ActiveDocument

# --- 2: plain text rewrite (single run, no interior <w:tab/>, so a normal
#     Find/Replace reproduces the target run structure exactly) ----------
$d.Content.Find.Execute(
    "My throughput improved from 73 operations in 100 milliseconds to 88 operations.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "My throughput improved from .73 operations per millisecond to .88 operations.",
    2) | Out-Null

# --- 3: "operations per 100 millisecs: 73" -> "avg operations per
#     millisec: .73". This text's run is preceded by a sibling <w:tab/>
#     element inside the very same run; a plain Find/Replace would fold
#     that <w:tab/> into the replaced <w:t> (turning it into a literal
#     tab character + xml:space="preserve"), which does not match the
#     target markup. Rewrite the whole paragraph's XML instead so the
#     <w:tab/> stays a distinct sibling, exactly like the source. ------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text -like "*operations per 100 millisecs: 73*") {
        $frag73 = (
            '<w:p>' +
              '<w:pPr>' +
                '<w:widowControl w:val="0"/>' +
                '<w:tabs>' +
                  '<w:tab w:val="left" w:pos="560"/>' +
                  '<w:tab w:val="left" w:pos="1120"/>' +
                  '<w:tab w:val="left" w:pos="1680"/>' +
                  '<w:tab w:val="left" w:pos="2240"/>' +
                  '<w:tab w:val="left" w:pos="2800"/>' +
                  '<w:tab w:val="left" w:pos="3360"/>' +
                  '<w:tab w:val="left" w:pos="3920"/>' +
                  '<w:tab w:val="left" w:pos="4480"/>' +
                  '<w:tab w:val="left" w:pos="5040"/>' +
                  '<w:tab w:val="left" w:pos="5600"/>' +
                  '<w:tab w:val="left" w:pos="6160"/>' +
                  '<w:tab w:val="left" w:pos="6720"/>' +
                '</w:tabs>' +
                '<w:autoSpaceDE w:val="0"/>' +
                '<w:autoSpaceDN w:val="0"/>' +
                '<w:adjustRightInd w:val="0"/>' +
                '<w:rPr>' +
                  '<w:rFonts w:ascii="Monaco" w:hAnsi="Monaco" w:cs="Monaco"/>' +
                  '<w:sz w:val="20"/>' +
                  '<w:szCs w:val="20"/>' +
                '</w:rPr>' +
              '</w:pPr>' +
              '<w:r>' +
                '<w:rPr>' +
                  '<w:rFonts w:ascii="Monaco" w:hAnsi="Monaco" w:cs="Monaco"/>' +
                  '<w:sz w:val="20"/>' +
                  '<w:szCs w:val="20"/>' +
                '</w:rPr>' +
                '<w:t xml:space="preserve">Running superscript </w:t>' +
              '</w:r>' +
              '<w:r>' +
                '<w:rPr>' +
                  '<w:rFonts w:ascii="Monaco" w:hAnsi="Monaco" w:cs="Monaco"/>' +
                  '<w:sz w:val="20"/>' +
                  '<w:szCs w:val="20"/>' +
                '</w:rPr>' +
                '<w:tab/>' +
              '</w:r>' +
              '<w:r>' +
                '<w:rPr>' +
                  '<w:rFonts w:ascii="Monaco" w:hAnsi="Monaco" w:cs="Monaco"/>' +
                  '<w:sz w:val="20"/>' +
                  '<w:szCs w:val="20"/>' +
                '</w:rPr>' +
                '<w:tab/>' +
              '</w:r>' +
              '<w:r>' +
                '<w:rPr>' +
                  '<w:rFonts w:ascii="Monaco" w:hAnsi="Monaco" w:cs="Monaco"/>' +
                  '<w:sz w:val="20"/>' +
                  '<w:szCs w:val="20"/>' +
                '</w:rPr>' +
                '<w:tab/>' +
              '</w:r>' +
              '<w:r>' +
                '<w:rPr>' +
                  '<w:rFonts w:ascii="Monaco" w:hAnsi="Monaco" w:cs="Monaco"/>' +
                  '<w:sz w:val="20"/>' +
                  '<w:szCs w:val="20"/>' +
                '</w:rPr>' +
                '<w:tab/>' +
                '<w:t>avg operations per millisec: .73</w:t>' +
              '</w:r>' +
            '</w:p>'
        )
        $cand.Range.InsertXML($frag73)
        break
    }
}

# --- 1: drop the "Also, for my script ..." paragraph plus the blank
#     paragraph right after it -------------------------------------------
$pAlso = $d.Paragraphs(6)
$pBlankAfter = $d.Paragraphs(7)
$toDelete = $d.Range($pAlso.Range.Start, $pBlankAfter.Range.End)
$toDelete.Delete()

# --- 4 & 5b: fix up the paragraph properties / runs on the last two
#     paragraphs of the document -----------------------------------------
$count = $d.Paragraphs.Count
$pSuperscript = $d.Paragraphs($count - 1)
$pFinal = $d.Paragraphs($count)

$superscriptXml = (
    '<w:p>' +
      '<w:pPr>' +
        '<w:widowControl w:val="0"/>' +
        '<w:autoSpaceDE w:val="0"/>' +
        '<w:autoSpaceDN w:val="0"/>' +
        '<w:adjustRightInd w:val="0"/>' +
        '<w:rPr>' +
          '<w:rFonts w:ascii="Monaco" w:hAnsi="Monaco" w:cs="Monaco"/>' +
          '<w:sz w:val="20"/>' +
          '<w:szCs w:val="20"/>' +
        '</w:rPr>' +
      '</w:pPr>' +
      '<w:r>' +
        '<w:rPr>' +
          '<w:rFonts w:ascii="Monaco" w:hAnsi="Monaco" w:cs="Monaco"/>' +
          '<w:sz w:val="20"/>' +
          '<w:szCs w:val="20"/>' +
        '</w:rPr>' +
        '<w:t xml:space="preserve">Running superscript </w:t>' +
      '</w:r>' +
      '<w:r>' +
        '<w:rPr>' +
          '<w:rFonts w:ascii="Monaco" w:hAnsi="Monaco" w:cs="Monaco"/>' +
          '<w:sz w:val="20"/>' +
          '<w:szCs w:val="20"/>' +
        '</w:rPr>' +
        '<w:tab/>' +
      '</w:r>' +
      '<w:r>' +
        '<w:rPr>' +
          '<w:rFonts w:ascii="Monaco" w:hAnsi="Monaco" w:cs="Monaco"/>' +
          '<w:sz w:val="20"/>' +
          '<w:szCs w:val="20"/>' +
        '</w:rPr>' +
        '<w:tab/>' +
      '</w:r>' +
      '<w:r>' +
        '<w:rPr>' +
          '<w:rFonts w:ascii="Monaco" w:hAnsi="Monaco" w:cs="Monaco"/>' +
          '<w:sz w:val="20"/>' +
          '<w:szCs w:val="20"/>' +
        '</w:rPr>' +
        '<w:tab/>' +
        '<w:t xml:space="preserve"> avg operations per millisec: .88</w:t>' +
      '</w:r>' +
    '</w:p>'
)
$pSuperscript.Range.InsertXML($superscriptXml)

$finalXml = (
    '<w:p>' +
      '<w:pPr>' +
        '<w:widowControl w:val="0"/>' +
        '<w:autoSpaceDE w:val="0"/>' +
        '<w:autoSpaceDN w:val="0"/>' +
        '<w:adjustRightInd w:val="0"/>' +
      '</w:pPr>' +
    '</w:p>'
)
$pFinal.Range.InsertXML($finalXml)
